# Updates the cryptos list with refreshed price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.543.02"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.595.73"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "1.822.23"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "1.583.31"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.539"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "27.529.70"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("D33").Value = "1.362.21"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.85%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").Value = "1.732.37"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0497"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
